$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 14
$ws.Range("H17").Value = 1194.3043
$ws.Range("J17").Value = 1228.591
$ws.Range("L17").Value = 3685.773
$ws.Range("N17").Value = -4021.773
$ws.Range("H88").Value = 4597.1113
$ws.Range("J88").Value = 4553.5713
$ws.Range("L88").Value = 4553.5713
$ws.Range("N88").Value = -5365.5713
$ws.Range("H91").Value = 4597.1113
$ws.Range("J91").Value = 4553.5713
$ws.Range("L91").Value = 4553.5713
$ws.Range("N91").Value = -7361.5713
$ws.Range("H137").Value = 3255.2307
$ws.Range("I137").Value = 2431.6667
$ws.Range("K137").Value = 7295.000100000001
$ws.Range("M137").Value = -4745.000100000001
$ws.Range("H138").Value = 3454.3635
$ws.Range("I138").Value = 1020.1539
$ws.Range("K138").Value = 3060.4617
$ws.Range("M138").Value = 2079.5383

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 40000
$ws.Range("I7").Value = 40000
$ws.Range("K7").Value = 40000
$ws.Range("M7").Value = -39886
$ws.Range("H122").Value = 4436.4707
$ws.Range("I122").Value = 4422.857
$ws.Range("K122").Value = 13268.571
$ws.Range("M122").Value = -10818.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 66712
$ws.Range("J13").Value = 66712
$ws.Range("L13").Value = 66712
$ws.Range("N13").Value = -67048
$ws.Range("H22").Value = 445.16666
$ws.Range("I22").Value = 449.54544
$ws.Range("J22").Value = 397
$ws.Range("K22").Value = 449.54544
$ws.Range("L22").Value = 397
$ws.Range("M22").Value = -276.54544
$ws.Range("N22").Value = -743
$ws.Range("H107").Value = 2079.2964
$ws.Range("I107").Value = 2099.16
$ws.Range("K107").Value = 2099.16
$ws.Range("M107").Value = -179.1599999999999
$ws.Range("H117").Value = 233246.67
$ws.Range("J117").Value = 233246.67
$ws.Range("L117").Value = 233246.67
$ws.Range("N117").Value = -242424.67
$ws.Range("H134").Value = 4308.467
$ws.Range("I134").Value = 3687.0908
$ws.Range("K134").Value = 11061.2724
$ws.Range("M134").Value = -8526.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2110.0527
$ws.Range("I31").Value = 1161.8462
$ws.Range("K31").Value = 1161.8462
$ws.Range("M31").Value = -866.8462
$ws.Range("H34").Value = 2110.0527
$ws.Range("I34").Value = 1161.8462
$ws.Range("K34").Value = 1161.8462
$ws.Range("M34").Value = -959.8462
$ws.Range("H58").Value = 1944.3784
$ws.Range("I58").Value = 1632.7646
$ws.Range("K58").Value = 1632.7646
$ws.Range("M58").Value = -1429.7646
$ws.Range("H62").Value = 3778.5715
$ws.Range("I62").Value = 3778.5715
$ws.Range("K62").Value = 3778.5715
$ws.Range("M62").Value = -3154.5715
$ws.Range("H65").Value = 3778.5715
$ws.Range("I65").Value = 3778.5715
$ws.Range("K65").Value = 18892.8575
$ws.Range("M65").Value = -15772.8575
$ws.Range("H95").Value = 15133
$ws.Range("J95").Value = 15133
$ws.Range("L95").Value = 15133
$ws.Range("N95").Value = -20625
$ws.Range("H99").Value = 8276.666999999999
$ws.Range("J99").Value = 8566.333000000001
$ws.Range("L99").Value = 8566.333000000001
$ws.Range("N99").Value = -11562.333
$ws.Range("H126").Value = 8276.666999999999
$ws.Range("J126").Value = 8566.333000000001
$ws.Range("L126").Value = 25698.999
$ws.Range("N126").Value = -30638.999
$ws.Range("H136").Value = 1944.3784
$ws.Range("I136").Value = 1632.7646
$ws.Range("K136").Value = 4898.293799999999
$ws.Range("M136").Value = -2348.293799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 1514
$ws.Range("I35").Value = 919.8
$ws.Range("J35").Value = 2999.5
$ws.Range("K35").Value = 2759.4
$ws.Range("L35").Value = 8998.5
$ws.Range("M35").Value = -2471.4
$ws.Range("N35").Value = -9574.5
$ws.Range("H123").Value = 3285.5
$ws.Range("I123").Value = 1240.6666
$ws.Range("K123").Value = 3721.9998
$ws.Range("M123").Value = -1271.9998
$ws.Range("H129").Value = 3273.44
$ws.Range("I129").Value = 956.1875
$ws.Range("K129").Value = 2868.5625
$ws.Range("M129").Value = 2131.4375
$ws.Range("H138").Value = 1258319.1
$ws.Range("I138").Value = 2511638.2
$ws.Range("J138").Value = 5000
$ws.Range("K138").Value = 7534914.600000001
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = -7529774.600000001
$ws.Range("N138").Value = -25280
$ws.Range("H139").Value = 4136.6
$ws.Range("I139").Value = 3381.125
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 10143.375
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = -5003.375
$ws.Range("N139").Value = -25280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 82666.664
$ws.Range("J100").Value = 82666.664
$ws.Range("L100").Value = 82666.664
$ws.Range("N100").Value = -84830.664
$ws.Range("H113").Value = 6765.6665
$ws.Range("I113").Value = 14998
$ws.Range("J113").Value = 2649.5
$ws.Range("K113").Value = 14998
$ws.Range("L113").Value = 2649.5
$ws.Range("M113").Value = -12828
$ws.Range("N113").Value = -6989.5
$ws.Range("H118").Value = 52322.25
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 52322.25
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 52322.25
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -55636.25
$ws.Range("H122").Value = 4470.4546
$ws.Range("I122").Value = 4422.1875
$ws.Range("K122").Value = 13266.5625
$ws.Range("M122").Value = -10816.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 30000
$ws.Range("I13").Value = 30000
$ws.Range("K13").Value = 30000
$ws.Range("M13").Value = -29860
$ws.Range("H22").Value = 2914.3572
$ws.Range("I22").Value = 732.6667
$ws.Range("J22").Value = 3509.3635
$ws.Range("K22").Value = 732.6667
$ws.Range("L22").Value = 3509.3635
$ws.Range("M22").Value = -437.6667
$ws.Range("N22").Value = -4099.363499999999
$ws.Range("H27").Value = 2914.3572
$ws.Range("I27").Value = 732.6667
$ws.Range("J27").Value = 3509.3635
$ws.Range("K27").Value = 732.6667
$ws.Range("L27").Value = 3509.3635
$ws.Range("M27").Value = -625.6667
$ws.Range("N27").Value = -3723.3635
$ws.Range("H122").Value = 7761.3335
$ws.Range("I122").Value = 7005.567
$ws.Range("J122").Value = 11540.167
$ws.Range("K122").Value = 21016.701
$ws.Range("L122").Value = 34620.501
$ws.Range("M122").Value = -18566.701
$ws.Range("N122").Value = -39520.501
$ws.Range("H132").Value = 4615.2354
$ws.Range("I132").Value = 3209.8572
$ws.Range("J132").Value = 5599
$ws.Range("K132").Value = 9629.571599999999
$ws.Range("L132").Value = 16797
$ws.Range("M132").Value = -7099.571599999999
$ws.Range("N132").Value = -21857
$ws.Range("H136").Value = 3796.1667
$ws.Range("I136").Value = 3669.7368
$ws.Range("K136").Value = 11009.2104
$ws.Range("M136").Value = -8459.2104

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2852.2354
$ws.Range("I96").Value = 2553.4546
$ws.Range("K96").Value = 2553.4546
$ws.Range("M96").Value = -1180.4546
$ws.Range("H122").Value = 2835.2285
$ws.Range("I122").Value = 2930.3225
$ws.Range("K122").Value = 8790.967500000001
$ws.Range("M122").Value = -6340.967500000001
$ws.Range("H132").Value = 948.069
$ws.Range("I132").Value = 817.6
$ws.Range("J132").Value = 1763.5
$ws.Range("K132").Value = 2452.8
$ws.Range("L132").Value = 5290.5
$ws.Range("M132").Value = 77.19999999999982
$ws.Range("N132").Value = -10350.5
$ws.Range("H136").Value = 4176.4165
$ws.Range("J136").Value = 5724.8
$ws.Range("L136").Value = 17174.4
$ws.Range("N136").Value = -22274.4
